# The journal entry for "4 fevrier" needs to be added to the document.
#
# In the source file the paragraph that should hold the "4 fevrier"
# heading already exists, but is empty and only carries a stray
# paragraph-mark run property (<w:pPr><w:rPr><w:lang w:val="en-US"/>
# </w:rPr></w:pPr>) -- no visible text and no bare <w:pPr>/<w:rPr> in the
# target. Simply assigning .Range.Text to it leaves that leftover <w:pPr>
# behind, so instead the paragraph's OOXML is replaced outright (via
# Range.InsertXML, i.e. Word's "flat OPC" insertion mechanism) with a
# clean paragraph that just contains the heading text.
#
# A second, brand-new paragraph is then inserted right after it for the
# journal entry itself. That paragraph is built from two separate runs
# (mirroring the source, which breaks "...le design" / " de la mallette"
# into two <w:r> elements) again using InsertXML, so the run split is
# preserved exactly rather than Word silently merging the text into a
# single run.

$d = $word.ActiveDocument

$xmlHeading = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>4 février</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xmlEntry = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Finir l’intégration des 3 esp32 et commencer le design</w:t></w:r><w:r><w:t xml:space="preserve"> de la mallette</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Paragraph 10 = the empty paragraph right before the trailing blank
# paragraph / sectPr (see diff: it currently only has <w:pPr><w:rPr>
# <w:lang w:val="en-US"/></w:rPr></w:pPr> and no runs).
$pHeading = $d.Paragraphs.Item(10)
$pHeading.Range.InsertXML($xmlHeading)

# The paragraph that used to be the trailing blank one is now number 11
# (paragraph 10 holds "4 fevrier"). Collapse to the very start of it and
# insert the journal-entry paragraph right there, pushing the blank
# paragraph further down instead of merging into it.
$pBlank = $d.Paragraphs.Item(11)
$rAfter = $pBlank.Range
$rAfter.Collapse(1)
$rAfter.InsertXML($xmlEntry)
